# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed values, per the GitHub Actions data-refresh job.
#
# Note: several Price values (e.g. "1.001", "237.26") are plain text in the
# source data (not real numbers - the sheet uses a locale where "." can be a
# thousands separator, e.g. "30.272.44"), so a leading apostrophe is used to
# force Excel to keep them as text instead of auto-converting them to
# numbers, and ClearFormats() strips the resulting quote-prefix cell format
# so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.272.44"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.865.71"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'237.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4684"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.2868"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").Value = "'0.06554"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'22.30"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +13.66%  "
$ws.Range("D11").Value = "'0.07902"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "1.869.64"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'5.180"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'0.6821"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "'278.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "30.281.45"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'13.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +8.73%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'5.399"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'0.000007347"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "2.112.44"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.205"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "'168.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("D26").Value = "'9.295"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'19.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "'1.949"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").Value = "'0.09851"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").Value = "'4.399"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'1.483"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'4.074"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'0.04752"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").Value = "'1.147"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("D36").Value = "'0.7097"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("D37").Value = "'2.707"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").Value = "'2.621"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("D40").Value = "'76.97"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.69%  "
$ws.Range("D41").Value = "'6.316"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "'1.966"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'0.4193"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").Value = "'0.9995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'103.47"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "'962.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "'7.228"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "'9.309"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'34.32"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").Value = "'0.05646"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.32%  "
